# BOT; UPDATE DATA
# Appends one new day's worth of PCR/infection figures (2020-05-01 serial
# 43948) to the three data sheets ("all", "kobe", "other"), pushing the
# trailing footnote row(s) down by one row on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all": new row 20 (date 43948) inserted above the two footnote
# rows, which shift from 20/21 down to 21/22.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Rows("20:20").Insert()
$wsAll.Range("A20").Value = 43948
$wsAll.Range("B20").Value = 251
$wsAll.Range("C20").Value = 224
$wsAll.Range("D20").Value = 133
$wsAll.Range("E20").Value = 123
$wsAll.Range("F20").Value = 10
$wsAll.Range("G20").Value = 3
$wsAll.Range("H20").Value = 88
$wsAll.Range("B23").Select()

# ---------------------------------------------------------------------
# Sheet "kobe": new row 75 (date 43948) inserted above the footnote row,
# which shifts from 75 down to 76. B75 (daily new-case count) is left
# blank, matching the source update.
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Rows("75:75").Insert()
$wsKobe.Range("A75").Value = 43948
$wsKobe.Range("C75").Value = 1681
$wsKobe.Range("D75").Value = 0
$wsKobe.Range("E75").Value = 251
$wsKobe.Range("F75").Value = 128
$wsKobe.Range("G75").Value = 119
$wsKobe.Range("H75").Value = 9
$wsKobe.Range("I75").Value = 3
$wsKobe.Range("J75").Value = 82
$wsKobe.Range("J76").Select()

# ---------------------------------------------------------------------
# Sheet "other": new row 50 (date 43948) inserted above the footnote
# row, which shifts from 50 down to 51; the previously-blank row 51
# shifts down to a new blank row 52.
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows("50:50").Insert()
$wsOther.Range("A50").Value = 43948
$wsOther.Range("B50").Value = 0
$wsOther.Range("C50").Value = 11
$wsOther.Range("D50").Value = 5
$wsOther.Range("E50").Value = 4
$wsOther.Range("F50").Value = 1
$wsOther.Range("G50").Value = 0
$wsOther.Range("H50").Value = 6
$wsOther.Range("H54").Select()

Write-Output "done"
